# Update the FHIR ValueSet "Metadata" sheet to the next IG release:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date bumped to the new publication timestamp
#  - Publisher contact info replaced with the actual publisher name, and a
#    "Jurisdiction" row (United States of America) replaces the old
#    "Contact" rows (one of which was a literal duplicate).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old sheet had two identical "Contact" / "No display for ContactDetail"
# rows (rows 10 and 11). Remove the duplicate row 11 entirely so everything
# below it shifts up by one, shrinking the sheet from 15 to 14 rows.
$ws.Rows.Item(11).Delete()

# Bump the Version value.
$ws.Range("B3").Value = "6.0.0"

# Bump the Date value.
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has an actual value.
$ws.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row becomes "Jurisdiction".
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
